# Updates the cryptos price-table rows (columns B-E) that changed between
# crawls. Most rows only touch Price (D) and Volume(1h) (E); two pairs of
# rows (37/38 and 45/46) also had their Coin name (B) and Link (C) swapped
# because the underlying ranking order changed.
#
# Price values are stored as literal text in the sheet (e.g. "1.00",
# "0.960", "58.695.81" with thousand-grouping dots) rather than numbers.
# Excel's COM layer auto-converts a plain-decimal-looking string assigned to
# .Value into a real number, which would silently drop formatting like
# trailing zeros. For those cells ("DIsNumericText") we briefly force a Text
# number format before writing the value, then restore the cell's style to
# "Normal" so no residual formatting/style is left on the cell afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Row=2; D="59.383.61"; E="  -1.43%  " },
    @{ Row=3; D="3.013.80"; E="  +1.51%  " },
    @{ Row=4; D="0.997"; DIsNumericText=$true; E="  -0.33%  " },
    @{ Row=5; D="565.38"; DIsNumericText=$true; E="  -0.70%  " },
    @{ Row=6; D="134.98"; DIsNumericText=$true; E="  +8.83%  " },
    @{ Row=7; D="0.996"; DIsNumericText=$true; E="  -0.43%  " },
    @{ Row=8; D="3.009.02"; E="  +1.48%  " },
    @{ Row=9; D="0.518"; DIsNumericText=$true; E="  +4.17%  " },
    @{ Row=10; D="0.132"; DIsNumericText=$true; E="  +0.28%  " },
    @{ Row=11; D="4.95"; DIsNumericText=$true; E="  -3.00%  " },
    @{ Row=12; D="0.456"; DIsNumericText=$true; E="  +4.91%  " },
    @{ Row=13; D="0.0000230"; DIsNumericText=$true; E="  +3.76%  " },
    @{ Row=14; D="33.85"; DIsNumericText=$true; E="  +4.44%  " },
    @{ Row=15; E="  +2.32%  " },
    @{ Row=16; D="3.490.34"; E="  +0.97%  " },
    @{ Row=17; D="6.88"; DIsNumericText=$true; E="  +12.03%  " },
    @{ Row=18; D="2.993.80"; E="  +0.74%  " },
    @{ Row=19; D="58.798.07"; E="  -2.51%  " },
    @{ Row=20; D="430.13"; DIsNumericText=$true; E="  +1.06%  " },
    @{ Row=21; D="13.40"; DIsNumericText=$true; E="  +3.11%  " },
    @{ Row=22; D="0.696"; DIsNumericText=$true; E="  +5.77%  " },
    @{ Row=23; D="7.13"; DIsNumericText=$true; E="  +1.22%  " },
    @{ Row=24; D="13.32"; DIsNumericText=$true; E="  +4.65%  " },
    @{ Row=25; D="80.62"; DIsNumericText=$true; E="  +2.70%  " },
    @{ Row=26; E="  +0.09%  " },
    @{ Row=27; E="  -0.02%  " },
    @{ Row=28; D="2.54"; DIsNumericText=$true; E="  +1.89%  " },
    @{ Row=29; D="7.83"; DIsNumericText=$true; E="  +10.06%  " },
    @{ Row=30; D="2.04"; DIsNumericText=$true; E="  +9.62%  " },
    @{ Row=31; D="6.41"; DIsNumericText=$true; E="  +6.03%  " },
    @{ Row=32; D="0.107"; DIsNumericText=$true; E="  +16.75%  " },
    @{ Row=33; D="25.63"; DIsNumericText=$true; E="  +2.08%  " },
    @{ Row=34; D="2.21"; DIsNumericText=$true; E="  -0.84%  " },
    @{ Row=35; D="5.76"; DIsNumericText=$true; E="  +4.12%  " },
    @{ Row=36; D="0.960"; DIsNumericText=$true; E="  +1.95%  " },
    @{ Row=37; B="PEPE"; C="https://coinranking.com/coin/03WI8NQPF+pepe-pepe"; D="0.0₃0713"; E="  +10.12%  " },
    @{ Row=38; B="OKB"; C="https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"; D="48.95"; DIsNumericText=$true; E="  -0.57%  " },
    @{ Row=39; D="8.50"; DIsNumericText=$true; E="  +8.52%  " },
    @{ Row=40; D="2.69"; DIsNumericText=$true; E="  +13.54%  " },
    @{ Row=41; D="0.112"; DIsNumericText=$true; E="  +2.55%  " },
    @{ Row=42; D="389.07"; DIsNumericText=$true; E="  +3.49%  " },
    @{ Row=43; D="0.0356"; DIsNumericText=$true; E="  +0.64%  " },
    @{ Row=44; D="2.686.58"; E="  +2.17%  " },
    @{ Row=45; B="TheGraph"; C="https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"; D="0.247"; DIsNumericText=$true; E="  +5.86%  " },
    @{ Row=46; B="USDe"; C="https://coinranking.com/coin/exbfr2U-0+usde-usde"; D="0.999"; DIsNumericText=$true; E="  +0.01%  " },
    @{ Row=47; D="2.05"; DIsNumericText=$true; E="  +5.07%  " },
    @{ Row=48; D="121.22"; DIsNumericText=$true; E="  +1.66%  " },
    @{ Row=49; E="  +4.47%  " },
    @{ Row=50; D="24.16"; DIsNumericText=$true; E="  +4.22%  " },
    @{ Row=51; D="2.05"; DIsNumericText=$true; E="  +4.16%  " }
)

foreach ($change in $changes) {
    $row = $change.Row

    foreach ($col in @("B", "C", "D", "E")) {
        if (-not $change.ContainsKey($col)) { continue }

        $value = $change[$col]
        $cell = $ws.Range("$col$row")

        if ($col -eq "D" -and $change.ContainsKey("DIsNumericText")) {
            # Force text storage so the literal string (with its original
            # formatting) is preserved instead of becoming a numeric value.
            $cell.NumberFormat = "@"
            $cell.Value = $value
            $cell.Style = "Normal"
        } else {
            $cell.Value = $value
        }
    }
}
